$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds plain-looking numeric strings (e.g. "7.240",
# "1.212", "0.9998"). Assigning them straight to .Value would let Excel
# auto-convert them to real floating point numbers, which drops trailing
# zeros / changes their representation. Force the column to text first.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$dRows = @(
    2,
    3,
    4,
    5,
    6,
    7,
    8,
    9,
    10,
    11,
    12,
    13,
    14,
    15,
    16,
    17,
    18,
    19,
    20,
    22,
    23,
    24,
    25,
    26,
    28,
    29,
    30,
    31,
    32,
    33,
    34,
    35,
    38,
    39,
    40,
    41,
    42,
    43,
    44,
    45,
    46,
    47,
    48,
    49,
    50,
    51
)
$dVals = @(
    '23.390.79',
    '1.631.05',
    '0.9998',
    '1.002',
    '297.79',
    '0.3756',
    '50.03',
    '0.3504',
    '0.08025',
    '1.206',
    '0.9996',
    '21.89',
    '6.292',
    '7.240',
    '0.00001192',
    '1.618.26',
    '95.03',
    '0.06928',
    '6.661',
    '1.001',
    '12.30',
    '23.402.94',
    '2.430',
    '2.951',
    '151.65',
    '5.159',
    '131.56',
    '1.813.04',
    '6.776',
    '2.104',
    '11.16',
    '0.9662',
    '0.2420',
    '5.845',
    '0.06748',
    '12.77',
    '0.6802',
    '1.290',
    '15.34',
    '0.9998',
    '0.6311',
    '2.229',
    '3.887',
    '0.07645',
    '126.44',
    '1.212'
)

$eRows = @(
    2,
    3,
    4,
    5,
    6,
    7,
    8,
    9,
    10,
    11,
    12,
    13,
    14,
    15,
    16,
    17,
    18,
    19,
    20,
    21,
    22,
    23,
    24,
    25,
    26,
    27,
    28,
    29,
    30,
    31,
    32,
    33,
    34,
    35,
    36,
    37,
    38,
    39,
    40,
    41,
    42,
    43,
    44,
    45,
    46,
    47,
    48,
    49,
    50,
    51
)
$eVals = @(
    '  -1.31%  ',
    '  -1.57%  ',
    '  -0.20%  ',
    '  -0.11%  ',
    '  -1.91%  ',
    '  -1.43%  ',
    '  -2.28%  ',
    '  -3.20%  ',
    '  -2.18%  ',
    '  -3.02%  ',
    '  -0.19%  ',
    '  -3.39%  ',
    '  -3.39%  ',
    '  -2.70%  ',
    '  -3.34%  ',
    '  -2.33%  ',
    '  -2.52%  ',
    '  -0.96%  ',
    '  -2.26%  ',
    '  -2.22%  ',
    '  -0.02%  ',
    '  -4.61%  ',
    '  -1.25%  ',
    '  -3.14%  ',
    '  -3.35%  ',
    '  -2.50%  ',
    '  -1.01%  ',
    '  -0.64%  ',
    '  -2.06%  ',
    '  -1.54%  ',
    '  -2.92%  ',
    '  -4.11%  ',
    '  -4.97%  ',
    '  -9.37%  ',
    '  -5.00%  ',
    '  -1.04%  ',
    '  -4.15%  ',
    '  -4.59%  ',
    '  -4.48%  ',
    '  -2.05%  ',
    '  -3.39%  ',
    '  -3.13%  ',
    '  -4.02%  ',
    '  -0.18%  ',
    '  -3.45%  ',
    '  -3.84%  ',
    '  -2.27%  ',
    '  -3.80%  ',
    '  -1.46%  ',
    '  +1.44%  '
)

for ($i = 0; $i -lt $dRows.Count; $i++) {
    $ws.Cells.Item($dRows[$i], 4).Value = $dVals[$i]
}

for ($i = 0; $i -lt $eRows.Count; $i++) {
    $ws.Cells.Item($eRows[$i], 5).Value = $eVals[$i]
}

# Restore the default "General" number format / style on column D so the
# saved workbook keeps the original (unstyled) appearance of these cells.
$priceRange.NumberFormat = "General"
$priceRange.Style = "Normal"

